# Update Il11-Il6st.xlsx worksheet with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5,6,7 (the old "ECs" sending-cluster rows have been replaced;
# only 3 data rows - all with Sending cluster = MuSCs - remain).
$ws.Rows.Item(5).Resize(3).Delete() | Out-Null

# Row 2: MuSCs / Il11 / Il6st / ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Il11"
$ws.Range("C2").Value = "Il6st"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.15518
$ws.Range("H2").Value = 0.46554
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 20.854426
$ws.Range("N2").Value = 62.563278
$ws.Range("O2").Value = 0.1507164072139519
$ws.Range("P2").Value = 0.1507164072139519
$ws.Range("Q2").Value = 3.23618982668
$ws.Range("R2").Value = 29.12570844012
$ws.Range("S2").Value = 0.1507164072139519
$ws.Range("T2").Value = 0.1507164072139519

# Row 3: MuSCs / Il11 / Il6st / FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Il11"
$ws.Range("C3").Value = "Il6st"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.15518
$ws.Range("H3").Value = 0.46554
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 94.96115633333334
$ws.Range("N3").Value = 284.883469
$ws.Range("O3").Value = 0.6862909728343718
$ws.Range("P3").Value = 0.6862909728343718
$ws.Range("Q3").Value = 14.73607223980667
$ws.Range("R3").Value = 132.62465015826
$ws.Range("S3").Value = 0.6862909728343718
$ws.Range("T3").Value = 0.6862909728343718

# Row 4: MuSCs / Il11 / Il6st / MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Il11"
$ws.Range("C4").Value = "Il6st"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.15518
$ws.Range("H4").Value = 0.46554
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.553069
$ws.Range("N4").Value = 67.65920700000001
$ws.Range("O4").Value = 0.1629926199516763
$ws.Range("P4").Value = 0.1629926199516763
$ws.Range("Q4").Value = 3.499785247420001
$ws.Range("R4").Value = 31.49806722678001
$ws.Range("S4").Value = 0.1629926199516763
$ws.Range("T4").Value = 0.1629926199516763
